$d = $word.ActiveDocument

# --- Paragraph 1: mark "sarapastrosa" / "esta" as spelling errors and "ojala" as a grammar error ---
$p1xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="70F6756B" w14:textId="504C1D21" w:rsidR="00F55F7C" w:rsidRDefault="00F265AB" w:rsidP="007164D1"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Esto es una prueba para la lectura de archivos con Python, fecha veintidós de noviembre del dos mil veintidós. Estoy con la </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>sarapastrosa</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> ella </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>esta</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> haciendo las interfaces </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>ojala</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> lo haga rápido.</w:t></w:r><w:r w:rsidR="007164D1"><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r1 = $d.Paragraphs(1).Range
$r1.InsertXML($p1xml)

# --- Paragraph 5: "Esta es la historia..." -> split out "que," as its own run ---
$p5xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="6D57A528" w14:textId="03E00134" w:rsidR="007164D1" w:rsidRDefault="007164D1" w:rsidP="007164D1"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Esta es la historia de un hidalgo de la Mancha de unos 50 años de edad </w:t></w:r><w:r><w:t>que,</w:t></w:r><w:r><w:t xml:space="preserve"> tras leer muchos libros de caballería, un género popular en siglo XVI, decide disfrazarse de caballero andante y embarcarse en una serie de aventuras al lado de su viejo caballo Rocinante. Tiene como fin &quot;irse por todo el mundo con sus armas y caballo a buscar las aventuras y a ejercitarse en todo aquello que él había leído que los caballeros andantes se ejercitaban, deshaciendo todo género de agravio y poniéndose en ocasiones y peligro donde, acabándolos, cobrase eterno nombre y fama&quot; (Parte 1, Cap. 1).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r5 = $d.Paragraphs(5).Range
$r5.InsertXML($p5xml)

# --- Paragraph 10: "Es una obra renacentista..." -> add page break + "Además," as its own run ---
$p10xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="13066829" w14:textId="33F5480C" w:rsidR="007164D1" w:rsidRDefault="007164D1" w:rsidP="007164D1"><w:pPr><w:jc w:val="both"/></w:pPr><w:r><w:t xml:space="preserve">Es una obra renacentista por su humanismo, pero con rasgos barrocos (el desengaño, el ambiente teatral y las apariencias falsas). Si bien es una parodia de libros de caballería en la que predomina el diálogo, también hay lugar para las historias intercaladas, que llegan gracias a otros personajes con quienes don Quijote se encuentra en sus andanzas. Estas nuevas historias </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">se prestan al empleo de diferentes estilos narrativos: el relato pastoril, la novela sentimental, la novela picaresca y la novela italiana son algunos ejemplos. </w:t></w:r><w:r><w:t>Además,</w:t></w:r><w:r><w:t xml:space="preserve"> aparece la tradición popular en los cuentos y refranes (“sabiduría popular”) de Sancho Panza. También incluye formas poéticas, como viejos romances caballerescos, canciones y sonetos.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r10 = $d.Paragraphs(10).Range
$r10.InsertXML($p10xml)

Write-Output "done"
